$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.579.97"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "3.069.51"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'550.40"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "'140.50"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.067.57"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'6.54"
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "'35.01"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "3.564.34"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").Value = "63.562.93"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "3.069.05"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("D21").Value = "'13.81"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'7.28"
$ws.Range("E23").Value = "  +4.84%  "
$ws.Range("D24").Value = "'81.15"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "'12.70"
$ws.Range("E25").Value = "  +5.77%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  +6.48%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'26.29"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "'2.47"
$ws.Range("E33").Value = "  +8.28%  "
$ws.Range("D34").Value = "'5.71"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").Value = "'55.55"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").Value = "'465.80"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").Value = "'0.0825"
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("D40").Value = "3.068.21"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "'2.57"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("D44").Value = "'28.16"
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'0.111"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'117.44"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0514"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "'2.09"
$ws.Range("E51").Value = "  +3.96%  "
